# Update cryptocurrency price/volume data on the active sheet.
# Values correspond to the refreshed data pulled by the scheduled GitHub
# Actions job on Fri Feb  3 08:42:12 UTC 2023.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "323.31";       E = "-1.95%" },
    @{ Row = 3;  D = "39.43";        E = "-1.34%" },
    @{ Row = 4;  D = "5.734";        E = "8.87%" },
    @{ Row = 5;  D = "0.07994";      E = "-1.36%" },
    @{ Row = 6;  D = "8.603";        E = "-0.40%" },
    @{ Row = 7;  D = "1.975";        E = "2.92%" },
    @{ Row = 8;  E = "-0.43%" },
    @{ Row = 9;  D = "0.9260";       E = "-1.01%" },
    @{ Row = 10; D = "0.1245";       E = "-6.23%" },
    @{ Row = 11; D = "0.1953";       E = "-0.87%" },
    @{ Row = 12; D = "8.713";        E = "24.87%" },
    @{ Row = 13; D = "0.09160";      E = "0.57%" },
    @{ Row = 14; D = "0.03652";      E = "2.86%" },
    @{ Row = 15; E = "9.44%" },
    @{ Row = 16; D = "0.001292";     E = "-1.95%" },
    @{ Row = 17; D = "0.006368";     E = "3.44%" },
    @{ Row = 18; D = "3.352";        E = "-0.36%" },
    @{ Row = 19; D = "4.528";        E = "-0.16%" },
    @{ Row = 20; D = "0.3537";       E = "0.59%" },
    @{ Row = 21; D = "0.1373";       E = "2.28%" },
    @{ Row = 22; D = "0.2450";       E = "-4.36%" },
    @{ Row = 23; D = "0.04409";      E = "-0.05%" },
    @{ Row = 24; D = "0.001263";     E = "3.33%" },
    @{ Row = 25; D = "0.004518";     E = "5.03%" },
    @{ Row = 26; D = "0.0001152";    E = "-3.17%" },
    @{ Row = 39; D = "0.02519";      E = "0.90%" },
    @{ Row = 40; D = "0.05413";      E = "4.47%" },
    @{ Row = 41; D = "0.007448";     E = "-3.56%" },
    @{ Row = 42; D = "0.009522";     E = "3.04%" },
    @{ Row = 43; D = "0.1404";       E = "-1.66%" },
    @{ Row = 44; D = "0.002124";     E = "-1.65%" },
    @{ Row = 45; D = "0.01070";      E = "4.57%" },
    @{ Row = 46; D = "0.00006778";   E = "1.75%" },
    @{ Row = 47; D = "0.00000000751";E = "0.16%" },
    @{ Row = 48; E = "-11.08%" },
    @{ Row = 49; D = "0.002294";     E = "-7.62%" },
    @{ Row = 50; D = "0.00002103";   E = "0.16%" },
    @{ Row = 51; D = "0.0002003";    E = "0.16%" }
)

# The sheet stores Price/Volume as plain text (e.g. "323.31", "-1.95%"),
# not numbers. Force text entry via NumberFormat "@" so Excel doesn't
# auto-convert them to numeric/percentage values, then restore the
# "Normal" style so no stray number-format style sticks to the cell
# (the original cells carry no explicit style).
foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        $cell = $ws.Range("D$($u.Row)")
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($u.ContainsKey("E")) {
        $cell = $ws.Range("E$($u.Row)")
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
        $cell.Style = "Normal"
    }
}
